# Update the binary "delivery sequence" assignment cells across the
# z{period},{zone} worksheets to reflect the re-solved model output.

$wb = $excel.ActiveWorkbook

# --- z1,1 ---
$ws = $wb.Worksheets.Item("z1,1")
$ws.Range("H5").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("G8").Value = 1

# --- z2,1 ---
$ws = $wb.Worksheets.Item("z2,1")
$ws.Range("F1").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("B6").Value = 1
$ws.Range("F10").Value = 0

# --- z1,2 ---
$ws = $wb.Worksheets.Item("z1,2")
$ws.Range("B6").Value = 0
$ws.Range("E6").Value = 0

# --- z1,3 ---
$ws = $wb.Worksheets.Item("z1,3")
$ws.Range("F1").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("F10").Value = 1

# --- z2,3 ---
$ws = $wb.Worksheets.Item("z2,3")
$ws.Range("H5").Value = 0
$ws.Range("C6").Value = 0
